# EPO output.xlsx — "improved code, fixed family members, view xml fucntion added"
#
# The sheet previously listed 3 patent records (rows 2-4) with columns:
#   B=Patent No. | C=Title | D=Family Id | E=Earliest priority | F=Publication date
#   | G=Assignee | H=Family members
#
# The updated sheet keeps only the first patent record (Oshkosh / US8139109B2),
# adds two new columns (Application Date, Inventors) and fixes the
# "Family members" list for that record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two extra patent rows (US5606609A, EP1417800B1) - only the first
# record (US8139109B2) survives.
$ws.Rows("3:4").Delete()

# Make room for the new "Application Date" column right after "Family Id"
# (old E:H -> new F:I).
$ws.Columns("E").Insert()

# Make room for the new "Inventors" column right after "Assignee"
# (old I = Family members -> new J).
$ws.Columns("I").Insert()

# --- Header row -------------------------------------------------------
$ws.Range("E1").Value = "Application Date"
$ws.Range("I1").Value = "Inventors"

# --- Data row (record 0) ----------------------------------------------
# Application date is a patent-office date code, stored as text (leading
# apostrophe forces Excel to keep it as text rather than converting it to
# a number), same convention already used for the publication date column.
$ws.Range("E2").Value = "'20061122"

# Inventors for the surviving record.
$ws.Range("I2").Value = "BROGGI ALBERTO, | SCHMIEDEL GARY, | YAKES CHRISTOPHER K"

# Corrected / de-duplicated family-members list for the surviving record.
$ws.Range("J2").Value = "US8947531B2| GB2473379A| GB2473379B| WO2008073518A2| CA2724324C| US2007291130A1| WO2009140514A3| WO2009140514A2| US8139109B2| CA2724324A1| US2012143430A1| GB201020969D0| US9420203B2| WO2008073518A3| US2009079839A1"
